$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the "Ají" dataset. It belongs
# chronologically above the existing row 663, so insert a fresh row there
# and push all subsequent records (old rows 663-720) down by one — this
# matches the diff exactly (old row N's data now lives in row N+1).
$ws.Rows(663).Insert()

$ws.Range("A663").Value = 3
$ws.Range("B663").Value = "Femacal de La Calera"
$ws.Range("C663").Value = "Coquimbo"
$ws.Range("D663").Value = 45106
$ws.Range("E663").Value = 5
$ws.Range("F663").Value = 100112021
$ws.Range("G663").Value = "Ají"
$ws.Range("H663").Value = "Inferno"
$ws.Range("I663").Value = "Primera"
$ws.Range("J663").Value = 73
$ws.Range("K663").Value = 10000
$ws.Range("L663").Value = 10500
$ws.Range("M663").Value = 10260
$ws.Range("N663").Value = "$/caja 10 kilos"
$ws.Range("O663").Value = "Región de Arica y Parinacota"
$ws.Range("P663").Value = 1026
$ws.Range("Q663").Value = 10
$ws.Range("R663").Value = "Hortaliza"
